{"js": "// Locate the last sentence of the conclusion (\"Zaklju\u010dak\") paragraph,\n// append two more sentences to that same paragraph, then add two more\n// new paragraphs (each consisting of two / one sentences) right after it,\n// all using the same body-text formatting (Times New Roman, 12pt /\n// sz 24 half-points, Serbian Latin language) that the rest of the\n// conclusion section already uses.\n\nconst body = context.document.body;\n\n// Find the existing final run of the conclusion paragraph so we can\n// anchor our insertions precisely after it.\nconst searchResults = body.search(\"propuste u okviru kategorije SSRF napada. \", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error(\"Could not find anchor text for the conclusion paragraph.\");\n}\n\nconst anchor = searchResults.items[0];\n\n// --- Append two sentences to the end of the existing conclusion paragraph ---\nconst afterAnchor = anchor.getRange(Word.RangeLocation.after);\nconst addedRun1 = afterAnchor.insertText(\n  \"Ova napadi demonstriraju da nekorektna obrada podataka koje unosi korisnik mo\u017ee da dovede do sigurnosnih propusta u organizaciji. \",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nconst afterRun1 = addedRun1.getRange(Word.RangeLocation.after);\nconst addedRun2 = afterRun1.insertText(\n  \"Posledice ovakvih napada su pristup internim servisima organizacije, zaobila\u017eenje firewall za\u0161tite i ostalih restrikcija i prikupljanje osetljivih podataka.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- Add a new paragraph with two more sentences ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst conclusionParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst newParagraph1 = conclusionParagraph.insertParagraph(\n  \"Za\u0161tita od korisni\u010dkog unosa se zasniva na filtriranju korisni\u010dkog unosa, odnosno odbacivanju zahtevanih URL koji se smatraju kao maliciozni. Glavni problem sa ovim je to \u0161to je sintaksa validno URL izuzetno slo\u017eena, tako da je lako prevariti URL parsere u svim programskim jezicima. \",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nconst newParagraph1Extra = newParagraph1.insertText(\n  \"Pored filtriranja postoje i druge tehnike, kao \u0161to je DNS rebinding, kojima se dodatno mo\u017ee zaobi\u0107i filtriranje URL.\",\n  Word.InsertLocation.end\n);\nawait context.sync();\n\n// --- Add the final concluding paragraph ---\nconst newParagraph2 = newParagraph1.insertParagraph(\n  \"Sa slo\u017eenijim strukturama organizacija i uvo\u0111enjem proxy servera ovaj tip napada je sve vi\u0161e zastupljen i postaje sve bitniji u razmatranju bezbednosti servera. \",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "# Add a conclusion to the \"Zaklju\u010dak\" (Conclusion) section:\n#  1) Append two more sentences to the existing final paragraph of the\n#     conclusion (the one ending in \"... propuste u okviru kategorije\n#     SSRF napada. \").\n#  2) Add a new paragraph with two sentences about input filtering and\n#     DNS rebinding.\n#  3) Add a final paragraph about growing relevance of this attack type.\n# All new text reuses the same formatting (Times New Roman, 12pt body\n# text, justified, first-line indent, Serbian Latin language) already\n# used throughout the rest of the document.\n\n$d = $word.ActiveDocument\n\n$wrapContinue = [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindContinue\n$replaceOne = [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceOne\n\n# --- 1) Extend the last sentence of the conclusion paragraph in place ---\n# Doing this through Find/Replace (rather than a plain InsertAfter at the\n# very end of the paragraph) makes the newly typed text merge into the\n# existing run and inherit its character formatting automatically.\n$findRange = $d.Content\n$findRange.Find.Execute(\n    \"propuste u okviru kategorije SSRF napada. \",\n    $false, $false, $false, $false, $false, $true, $wrapContinue, $false,\n    \"propuste u okviru kategorije SSRF napada. Ova napadi demonstriraju da nekorektna obrada podataka koje unosi korisnik mo\u017ee da dovede do sigurnosnih propusta u organizaciji. Posledice ovakvih napada su pristup internim servisima organizacije, zaobila\u017eenje firewall za\u0161tite i ostalih restrikcija i prikupljanje osetljivih podataka.\",\n    $replaceOne\n) | Out-Null\n\n# --- 2) Add a new paragraph right after the conclusion paragraph ---\n$conclusionParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$conclusionParagraph.Range.InsertParagraphAfter()\n$newParagraph1 = $d.Paragraphs.Item($d.Paragraphs.Count)\n$newParagraph1.Range.Text = \"Za\u0161tita od korisni\u010dkog unosa se zasniva na filtriranju korisni\u010dkog unosa, odnosno odbacivanju zahtevanih URL koji se smatraju kao maliciozni. Glavni problem sa ovim je to \u0161to je sintaksa validno URL izuzetno slo\u017eena, tako da je lako prevariti URL parsere u svim programskim jezicima. Pored filtriranja postoje i druge tehnike, kao \u0161to je DNS rebinding, kojima se dodatno mo\u017ee zaobi\u0107i filtriranje URL.\"\n\n# --- 3) Add the final paragraph of the conclusion ---\n$newParagraph1 = $d.Paragraphs.Item($d.Paragraphs.Count)\n$newParagraph1.Range.InsertParagraphAfter()\n$newParagraph2 = $d.Paragraphs.Item($d.Paragraphs.Count)\n$newParagraph2.Range.Text = \"Sa slo\u017eenijim strukturama organizacija i uvo\u0111enjem proxy servera ovaj tip napada je sve vi\u0161e zastupljen i postaje sve bitniji u razmatranju bezbednosti servera. \"\n"}
